# Applies the commit:
#  - "no longer use noise to generate terrain so removed the library"
#  - removes stray <w:proofErr/> spell-check markers around GeoJSON,
#    openGL, Imgui and MyGAL
#  - removes the "FastNoiseLite" library paragraph (and the extra blank
#    paragraph that followed it) since terrain generation no longer
#    depends on that library

$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'
$w14Ns = 'http://schemas.microsoft.com/office/word/2010/wordml'

$cr = [char]13

function Get-ParaByText($doc, [string]$text) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.TrimEnd($cr) -eq $text) {
            return $p
        }
    }
    return $null
}

# --- 1. "Description" paragraph: drop the spellStart/spellEnd proofErr
#        markers bracketing "GeoJSON" while keeping the surrounding runs
#        ("in " / "GeoJSON" / " format") intact and separate. -----------
$descPara = Get-ParaByText $d "Description: Randomly generate 2D maps and nations that occupy the land. Using files in GeoJSON format be able to save and load maps."
if ($descPara -eq $null) { throw "Could not find Description paragraph" }

$descXml = '<w:p xmlns:w="' + $wNs + '" xmlns:w14="' + $w14Ns + '" w14:paraId="6BE69D95" w14:textId="2ADA2635" w:rsidR="0044446F" w:rsidRDefault="0044446F" w:rsidP="0044446F">' +
    '<w:r><w:t>Description: Randomly generate 2D maps</w:t></w:r>' +
    '<w:r w:rsidR="00AF3318"><w:t xml:space="preserve"> and nations that occupy the land</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">. Using files </w:t></w:r>' +
    '<w:r w:rsidR="00566BFC"><w:t xml:space="preserve">in </w:t></w:r>' +
    '<w:r><w:t>GeoJSON</w:t></w:r>' +
    '<w:r w:rsidR="00566BFC"><w:t xml:space="preserve"> format</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> be able to save and load maps.</w:t></w:r>' +
    '</w:p>'
$descPara.Range.InsertXML($descXml)

# --- 2. "Libraries (...)" paragraph: proofErr removed AND the three
#        runs collapse into a single run. -------------------------------
$libPara = Get-ParaByText $d "Libraries (apart from all the openGL ones)"
if ($libPara -eq $null) { throw "Could not find Libraries paragraph" }

$libXml = '<w:p xmlns:w="' + $wNs + '" xmlns:w14="' + $w14Ns + '" w14:paraId="6BC15465" w14:textId="0652AF35" w:rsidR="00566BFC" w:rsidRDefault="00566BFC" w:rsidP="0044446F">' +
    '<w:r><w:t>Libraries (apart from all the openGL ones)</w:t></w:r>' +
    '</w:p>'
$libPara.Range.InsertXML($libXml)

# --- 3. "Dear Imgui." paragraph: proofErr removed AND the runs collapse
#        into a single run. ----------------------------------------------
$imguiPara = Get-ParaByText $d "Dear Imgui."
if ($imguiPara -eq $null) { throw "Could not find Dear Imgui paragraph" }

$imguiXml = '<w:p xmlns:w="' + $wNs + '" xmlns:w14="' + $w14Ns + '" w14:paraId="0057E877" w14:textId="30A45277" w:rsidR="00566BFC" w:rsidRDefault="00566BFC" w:rsidP="0044446F">' +
    '<w:r><w:t>Dear Imgui.</w:t></w:r>' +
    '</w:p>'
$imguiPara.Range.InsertXML($imguiXml)

# --- 4. "MyGAL" paragraph: just drop the surrounding proofErr markers,
#        run/text stays the same. ----------------------------------------
$myGalPara = Get-ParaByText $d "MyGAL"
if ($myGalPara -eq $null) { throw "Could not find MyGAL paragraph" }

$myGalXml = '<w:p xmlns:w="' + $wNs + '" xmlns:w14="' + $w14Ns + '" w14:paraId="2AFC643E" w14:textId="16ADB003" w:rsidR="00940180" w:rsidRDefault="00940180" w:rsidP="0044446F">' +
    '<w:r><w:t>MyGAL</w:t></w:r>' +
    '</w:p>'
$myGalPara.Range.InsertXML($myGalXml)

# --- 5. Remove the "FastNoiseLite" paragraph entirely together with the
#        blank paragraph immediately following it (terrain no longer
#        generated with noise, so the library + its trailing spacer
#        paragraph both go). ---------------------------------------------
$fnlPara = Get-ParaByText $d "FastNoiseLite"
if ($fnlPara -eq $null) { throw "Could not find FastNoiseLite paragraph" }

$rangeStart = $fnlPara.Range.Start
$nextPara = $fnlPara.Next()
$rangeEnd = $fnlPara.Range.End
if ($nextPara -ne $null -and $nextPara.Range.Text.TrimEnd($cr) -eq "") {
    $rangeEnd = $nextPara.Range.End
}
$d.Range($rangeStart, $rangeEnd).Delete()
